# Fruta / hortaliza, semanal
# Insert a new weekly data row into the table at row 174 (pushing the
# existing rows 174-177 down to 175-178) and populate it with the new
# week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 174; this shifts rows
# 174-177 down to 175-178, keeping all their original values intact.
$ws.Rows.Item(174).Insert()

# Populate the newly inserted row 174 with this week's record.
$ws.Range("A174").Value = 4
$ws.Range("B174").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C174").Value = "Los Lagos"
$ws.Range("D174").Value = 44448
$ws.Range("E174").Value = 10
$ws.Range("F174").Value = 100114014
$ws.Range("G174").Value = "Betarraga"
$ws.Range("H174").Value = "Sin especificar"
$ws.Range("I174").Value = "Primera"
$ws.Range("J174").Value = 500
$ws.Range("K174").Value = 1000
$ws.Range("L174").Value = 1000
$ws.Range("M174").Value = 1000
$ws.Range("N174").Value = '$/paquete 5 unidades'
$ws.Range("O174").Value = "Región del Maule"
$ws.Range("P174").Value = 200
$ws.Range("Q174").Value = 5
$ws.Range("R174").Value = "Hortaliza"
